$d = $word.ActiveDocument

# Raw OOXML for the 4 new paragraphs appended after the "Loi ich - Kho khan"
# paragraph (commit: "update project: ket luan, link source").
$newParagraphsXml = @'
<w:p><w:pPr><w:pStyle w:val="Heading1"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:bookmarkStart w:id="14" w:name="_Toc148800175"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs w:val="0"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr><w:t>Chú thích:</w:t></w:r><w:bookmarkEnd w:id="14"/></w:p>
<w:p><w:pPr><w:rPr><w:rFonts w:cs="Times New Roman"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b/><w:bCs/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr><w:t xml:space="preserve">[1] Same origin policy: </w:t></w:r><w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr><w:t>Chính sách cùng nguồn gốc là một cơ chế an ninh trong trình duyệt web để ngăn chặn các trang web hoặc tài nguyên từ nguồn không cùng nguồn gốc truy cập và tương tác với các dữ liệu từ nguồn khác. Nguồn gốc ở đây thường được xác định bằng cách so sánh giao thức, tên miền và cổng của hai trang web.</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:rFonts w:cs="Times New Roman"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr><w:t>Ví dụ: Giả sử bạn đang duyệt một trang web A từ miền "https://www.example.com". Trang web A muốn tải dữ liệu từ một trang web B có nguồn gốc từ "https://api.example.com". Nếu trang web B không chia sẻ cùng nguồn gốc với trang web A, thì trình duyệt sẽ áp dụng chính sách cùng nguồn gốc và ngăn chặn trang web A truy cập dữ liệu từ trang web B. Điều này là để đảm bảo rằng trang web A không thể lợi dụng quyền truy cập vào dữ liệu nhạy cảm trên trang web B mà người dùng không được phép truy cập.</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="vi-VN"/></w:rPr><w:t>Trong ngữ cảnh của Selenium, vấn đề xuất phát khi người kiểm thử cố gắng kiểm tra ứng dụng web trên một trang web khác nguồn gốc. Ví dụ, nếu bạn đang sử dụng Selenium để kiểm tra một trang web từ miền "https://www.example.com", nhưng Selenium được thực thi từ một miền khác, chẳng hạn "https://test.example". Chính sách cùng nguồn gốc sẽ ngăn chặn Selenium truy cập và tương tác với các thành phần của trang web từ miền "https://www.example.com".</w:t></w:r></w:p>
'@

$endRange = $d.Range($d.Content.End, $d.Content.End)
$endRange.InsertXML($newParagraphsXml)
